$wb = $excel.ActiveWorkbook

# --- 1. Rename the "Include from SNOMED CT[ N]" sheets to "Include #N" ---
$includeNames = @(
    "Include from SNOMED CT",
    "Include from SNOMED CT 2",
    "Include from SNOMED CT 3",
    "Include from SNOMED CT 4",
    "Include from SNOMED CT 5",
    "Include from SNOMED CT 6",
    "Include from SNOMED CT 7",
    "Include from SNOMED CT 8",
    "Include from SNOMED CT 9"
)

for ($i = 0; $i -lt $includeNames.Length; $i++) {
    $sheet = $wb.Worksheets.Item($includeNames[$i])
    $sheet.Name = "Include #$i"
}

# --- 2. Update the Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: pythia -> cicada
$meta.Cells.Item(2, 2).Value = "http://fhirfli.dev/fhir/ig/cicada/ValueSet/vaccine-medication-codes-snomed"

# Date: refreshed IG build timestamp
$meta.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# Insert a new "Jurisdiction" row right after "Contact" (row 10), before "Description" (row 11)
$meta.Rows.Item(11).Insert()

$meta.Cells.Item(11, 1).Value = "Jurisdiction"
$meta.Cells.Item(11, 2).Value = ""

# Match the formatting of the surrounding property rows
$meta.Range("A11:B11").Style = $meta.Range("A10:B10").Style
